# Cash Management test-data workbook: the author stripped the hard-coded
# Oracle Cloud URL / username / password out of the "Input_Value" sheet
# before re-uploading it, leaving those three cells blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Clear the URL, UserName and Password values that used to live in M2:O2
# (columns M/N/O are headed "URL" / "UserName" / "Password" in row 1).
$ws.Range("M2:O2").ClearContents()

# Leave the same range selected/active, matching the on-screen state the
# workbook was saved with.
$ws.Activate()
$ws.Range("M2:O2").Select()
